$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) slide1.xml: merge adjacent same-format runs into single runs
#    ("Etude de l'existant " + " " -> "Etude de l'existant  ")
#    ("Glasir : cahier des " + "charges`t" -> "Glasir : cahier des charges`t")
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$body1 = $s1.Shapes.Item(2).TextFrame.TextRange

$merge1 = $body1.Characters(43, 21)
$merge1.Text = "Etude de l’existant  "

$merge2 = $body1.Characters(72, 28)
$merge2.Text = "Glasir : cahier des charges`t"

# ---------------------------------------------------------------------------
# 2) Add a new slide (id 260) right after slide 1, using the "Titre et
#    contenu" layout (same layout as slide 1) and fill in the content
#    placeholder text.
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Add(2, 2)

$content = $s2.Shapes.Item(2).TextFrame.TextRange
$content.Text = " 16 et 18 mai 2012`r`r`r`rLe métro rennais paralysé`r`r"
$content.Font.Size = 48
$content.ParagraphFormat.Alignment = 2

# Split the first paragraph back into two runs: " " and "16 et 18 mai 2012"
$firstSpace = $content.Characters(1, 1)
$firstSpace.Text = " "
